$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B19").Value = "test"
